# Edit: "Added call to help Transcrypt gain mindshare (7)"
#
# The deck has a single slide with one large text box (shape 1) that
# contains the "spread the word" copy. Two paragraphs need their run
# structure changed:
#
#   1. The "Help it gain even more mindshare..." paragraph currently
#      is split across two runs ("H" / "elp it gain...") that must be
#      merged into a single run.
#
#   2. The "- Write about it, use it, let me know what you use it for,"
#      paragraph is currently a single run and must be split into four
#      runs, changing "me" to "us" along the way:
#        "- Write about it, use it, let " / "us" / " " / "know what you use it for,"
#
# Note: TextRange.Text round-trips curly quotes as plain ASCII apostrophes
# when *read back*, even though the underlying XML stores the real
# U+2019 glyph, so paragraph matching below compares against the plain
# apostrophe form while the text that gets *written* back uses the
# proper curly apostrophe (so the saved OOXML keeps matching byte-for-byte).

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(1)
$tf = $sh.TextFrame
$tr = $tf.TextRange

$rsquo = [char]8217

# --- Paragraph: "Help it gain even more mindshare, ..." -------------------
# Force the paragraph down to a single run. Re-writing a paragraph's Text
# with content that is byte-identical to what's already there leaves the
# existing run split alone, so first overwrite it with throwaway text
# (collapsing it to one run), then set it to the real final text.

$count = $tr.Paragraphs().Count
for ($i = 1; $i -le $count; $i++) {
    $para = $tr.Paragraphs($i, 1)
    $t = $para.Text.TrimEnd([char]13)
    if ($t -eq "Help it gain even more mindshare, because that's what'll convince companies to rely on it.") {
        $para.Text = "TEMP_PLACEHOLDER_COLLAPSE_RUNS"
        $para2 = $tr.Paragraphs($i, 1)
        $para2.Text = "Help it gain even more mindshare, because that" + $rsquo + "s what" + $rsquo + "ll convince companies to rely on it."
        break
    }
}

# --- Paragraph: "- Write about it, use it, let me know what you use it for," ---
# Split the single run into four runs, changing "me" to "us":
#   "- Write about it, use it, let " | "us" | " " | "know what you use it for,"

$count = $tr.Paragraphs().Count
for ($i = 1; $i -le $count; $i++) {
    $para = $tr.Paragraphs($i, 1)
    $t = $para.Text.TrimEnd([char]13)
    if ($t -eq "- Write about it, use it, let me know what you use it for,") {
        $st = $para.Start

        $part1 = $tr.Characters($st, 30)
        $part1.Text = "- Write about it, use it, let "

        $part2 = $tr.Characters($st + 30, 2)
        $part2.Text = "us"

        $part3 = $tr.Characters($st + 32, 1)
        $part3.Text = " "

        $part4 = $tr.Characters($st + 33, 25)
        $part4.Text = "know what you use it for,"

        break
    }
}

# --- Restore the autofit height -------------------------------------------
# The shape has <a:spAutoFit/>, so editing its text makes this host
# recompute the rendered height, same as real PowerPoint would on an
# autofit text box. The source deck's cached extent (cy="6524863" EMU,
# i.e. 513.76874... pt) predates this edit though, and the diff we're
# replaying didn't touch that cached value, so put it back afterwards.
$sh.Height = 513.7688

